# Updates the "cryptos" price list with refreshed prices / 1h volume figures,
# and fixes the ordering of the Polkadot / WrappedEther rows (12 and 13).
# Price (column D) values are kept as plain text (matching the source data,
# which stores prices like "26.124.79" / "0.2620" as text, not numbers), so
# NumberFormat is forced to Text ("@") before assignment and then the cell
# style is reset back to Normal so no other formatting is changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.124.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5311"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.500"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.662.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5482"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.127.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.560"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.030"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "140.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1247"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.283"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.434"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05960"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.281"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.518"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.245"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.562"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9528"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.34%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -4.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01614"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8465"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.016.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.800.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈102"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.40%  "
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.476"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.773"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.10%  "
